# Add force calibrate function
#
# Inserts a new Modbus register row (220, MB_FORCE_CALIBRATE_ENCODER,)
# directly after register 219 (MB_GOTO_SPEED_SETPOINT,) on the "Sheet1"
# map, pushing every following row down by one (register 9007 /
# MB_MAX_CURRENT_LIMIT_OUTWARD, now lands on row 73 instead of the old
# trailing blank row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the old row 51 (was register 299 /
# MB_EXTENSION,) — this shifts rows 51..72 down to 52..73 and keeps every
# existing cell (value + style) intact.
$ws.Rows("51:51").Insert()

# Populate the newly inserted row with the new register definition.
$ws.Range("A51").Value = 220
$ws.Range("B51").Value = "MB_FORCE_CALIBRATE_ENCODER,"
$ws.Range("C51").Value = "Write 0xA0A0 to force encoder to calibrate to zero in current position"
$ws.Range("D51").Value = "W"

# Match the author's final selection / scroll position.
$ws.Range("E51").Select()
$excel.ActiveWindow.ScrollRow = 35
